$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto quotes (price + 1h volume change).
#
# Every Price (column D) write below is prefixed with a literal leading
# apostrophe so Excel stores the cell as quote-prefixed TEXT instead of
# silently coercing it to a Number. The sheet intentionally keeps these as
# text (many entries use a non-numeric thousand-separator style like
# "27.080.25", and several have significant trailing zeros, e.g. "1.000" /
# "0.07570"), so a plain numeric write would corrupt the displayed value.

$ws.Cells.Item(2, 4).Value = '''27.070.63'
$ws.Cells.Item(2, 5).Value = '  -2.35%  '
$ws.Cells.Item(3, 4).Value = '''1.865.44'
$ws.Cells.Item(3, 5).Value = '  -2.07%  '
$ws.Cells.Item(4, 4).Value = '''1.001'
$ws.Cells.Item(4, 5).Value = '  +0.15%  '
$ws.Cells.Item(5, 4).Value = '''306.54'
$ws.Cells.Item(5, 5).Value = '  -1.96%  '
$ws.Cells.Item(6, 5).Value = '  +0.06%  '
$ws.Cells.Item(7, 4).Value = '''0.5115'
$ws.Cells.Item(7, 5).Value = '  -1.77%  '
$ws.Cells.Item(8, 4).Value = '''0.3747'
$ws.Cells.Item(8, 5).Value = '  -0.82%  '
$ws.Cells.Item(9, 4).Value = '''0.07156'
$ws.Cells.Item(9, 5).Value = '  -1.14%  '
$ws.Cells.Item(10, 4).Value = '''0.8884'
$ws.Cells.Item(10, 5).Value = '  -1.29%  '
$ws.Cells.Item(11, 4).Value = '''20.68'
$ws.Cells.Item(11, 5).Value = '  -2.81%  '
$ws.Cells.Item(12, 4).Value = '''0.07570'
$ws.Cells.Item(12, 5).Value = '  -0.87%  '
$ws.Cells.Item(13, 4).Value = '''1.853.04'
$ws.Cells.Item(13, 5).Value = '  -2.57%  '
$ws.Cells.Item(14, 5).Value = '  -2.60%  '
$ws.Cells.Item(15, 4).Value = '''89.47'
$ws.Cells.Item(15, 5).Value = '  -2.85%  '
$ws.Cells.Item(16, 4).Value = '''1.001'
$ws.Cells.Item(16, 5).Value = '  +0.17%  '
$ws.Cells.Item(17, 4).Value = '''0.000008447'
$ws.Cells.Item(17, 5).Value = '  -2.95%  '
$ws.Cells.Item(18, 5).Value = '  -2.65%  '
$ws.Cells.Item(19, 4).Value = '''1.000'
$ws.Cells.Item(19, 5).Value = '  +0.06%  '
$ws.Cells.Item(20, 4).Value = '''27.128.26'
$ws.Cells.Item(20, 5).Value = '  -2.26%  '
$ws.Cells.Item(21, 4).Value = '''5.036'
$ws.Cells.Item(21, 5).Value = '  -2.10%  '
$ws.Cells.Item(22, 4).Value = '''2.095.47'
$ws.Cells.Item(22, 5).Value = '  -2.02%  '
$ws.Cells.Item(23, 4).Value = '''10.54'
$ws.Cells.Item(23, 5).Value = '  -2.72%  '
$ws.Cells.Item(24, 4).Value = '''6.456'
$ws.Cells.Item(25, 4).Value = '''1.844'
$ws.Cells.Item(25, 5).Value = '  -1.88%  '
$ws.Cells.Item(26, 4).Value = '''147.88'
$ws.Cells.Item(26, 5).Value = '  -3.58%  '
$ws.Cells.Item(27, 4).Value = '''17.98'
$ws.Cells.Item(27, 5).Value = '  -1.84%  '
$ws.Cells.Item(28, 4).Value = '''2.114'
$ws.Cells.Item(28, 5).Value = '  -2.48%  '
$ws.Cells.Item(29, 4).Value = '''112.88'
$ws.Cells.Item(29, 5).Value = '  -1.30%  '
$ws.Cells.Item(30, 4).Value = '''4.666'
$ws.Cells.Item(30, 5).Value = '  -4.09%  '
$ws.Cells.Item(31, 4).Value = '''4.705'
$ws.Cells.Item(31, 5).Value = '  -3.10%  '
$ws.Cells.Item(32, 4).Value = '''0.09094'
$ws.Cells.Item(32, 5).Value = '  +1.24%  '
$ws.Cells.Item(33, 5).Value = '  -2.92%  '
$ws.Cells.Item(34, 4).Value = '''3.052'
$ws.Cells.Item(34, 5).Value = '  -3.92%  '
$ws.Cells.Item(35, 4).Value = '''1.155'
$ws.Cells.Item(35, 5).Value = '  -6.07%  '
$ws.Cells.Item(36, 4).Value = '''0.7270'
$ws.Cells.Item(36, 5).Value = '  -6.00%  '
$ws.Cells.Item(37, 5).Value = '  -2.11%  '
$ws.Cells.Item(38, 4).Value = '''2.495'
$ws.Cells.Item(38, 5).Value = '  -5.30%  '
$ws.Cells.Item(39, 4).Value = '''3.043'
$ws.Cells.Item(39, 5).Value = '  -0.73%  '
$ws.Cells.Item(40, 5).Value = '  -1.67%  '
$ws.Cells.Item(41, 4).Value = '''0.5335'
$ws.Cells.Item(41, 5).Value = '  -3.29%  '
$ws.Cells.Item(42, 4).Value = '''6.570'
$ws.Cells.Item(42, 5).Value = '  -1.54%  '
$ws.Cells.Item(43, 4).Value = '''116.98'
$ws.Cells.Item(43, 5).Value = '  +2.13%  '
$ws.Cells.Item(44, 4).Value = '''8.270'
$ws.Cells.Item(44, 5).Value = '  -2.87%  '
$ws.Cells.Item(45, 4).Value = '''0.1471'
$ws.Cells.Item(45, 5).Value = '  -2.59%  '
$ws.Cells.Item(48, 4).Value = '''9.992'
$ws.Cells.Item(48, 5).Value = '  -3.72%  '
$ws.Cells.Item(49, 4).Value = '''1.569'
$ws.Cells.Item(49, 5).Value = '  -2.82%  '
$ws.Cells.Item(50, 4).Value = '''36.55'
$ws.Cells.Item(50, 5).Value = '  -0.69%  '
$ws.Cells.Item(51, 4).Value = '''63.92'
$ws.Cells.Item(51, 5).Value = '  -4.27%  '

# Rows 46/47 swap rank order: Decentraland now ranks above PaxDollar.
$ws.Cells.Item(46, 2).Value = 'Decentraland'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(46, 4).Value = '''0.4635'
$ws.Cells.Item(46, 5).Value = '  -3.57%  '

$ws.Cells.Item(47, 2).Value = 'PaxDollar'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(47, 4).Value = '''0.9999'
$ws.Cells.Item(47, 5).Value = '  +0.05%  '

